# Importer: prevent product to be linked to variation child
# Insert a new test row right before the former row 46 ("Color/Black" /
# "Size/XS" variation) that exercises a product whose parent SKU would
# link it to a variation child. This pushes the previously existing rows
# 46-50 down to 47-51, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 46, shifting existing rows 46-50 down
# to 47-51.
$ws.Rows(46).Insert()

# Populate the new row 46 with the new test case data.
$ws.Range("B46").Value = 42
$ws.Range("C46").Value = 28
$ws.Range("D46").Value = "This tries to link variation to child"
$ws.Range("F46").Value = "Color/Black"
$ws.Range("G46").Value = "Size/XS"
$ws.Range("I46").Value = 12
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = "Test Category"
$ws.Range("L46").Value = "Test Category"
$ws.Range("M46").Value = "shirt1.jpeg"
$ws.Range("N46").Value = "shirt2.jpeg,shirt3.jpeg"

# Match the cursor position left behind by the authoring session.
$ws.Range("B52").Select()
